$wb = $excel.ActiveWorkbook

# The status text "Ready for handoff" for the 06e6f787 file is one shared
# string referenced by all four of these cells; the handback transform
# failure flips that single piece of text everywhere it is shown.
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: populate Error Detail (column P) for the 06e6f787 row ---
$zhcn.Range("P3").Value = "Handback file name: ct0v4vwn.sfp is different with handoff file name: 06e6f787-e5cf-48d3-9e97-84bd3f7dbb6e.4dbde9024bd557db508680c3db837c542d9e572a.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: populate Error Detail (column P) for the 06e6f787 row ---
$dede.Range("P3").Value = "Handback file name: ct0v4vwn.sfp is different with handoff file name: 06e6f787-e5cf-48d3-9e97-84bd3f7dbb6e.4dbde9024bd557db508680c3db837c542d9e572a.de-de."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
